# Rename the worksheet from "Sheet1" to "majors".
# (This also updates the workbook-scoped _xlnm._FilterDatabase defined name,
# whose formula references the sheet by name.)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "majors"

# Reflect the author's active cell selection on the sheet (cell A87).
$ws.Range("A87").Select()
